$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.266970157623291
$ws.Range("B1").Value = 1.748167037963867
$ws.Range("C1").Value = 4.080344676971436
$ws.Range("D1").Value = 3.331408500671387
$ws.Range("E1").Value = 1.150517821311951
